$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The worker table (rows 16-25, 10 employees) gets replaced with an updated
# set of 12 employees (some kept, some new, reordered). We insert two new
# rows right after the current last data row (25) so the table grows from
# 10 to 12 rows, then fix up formatting, then write all 12 rows of data.
# ---------------------------------------------------------------------------

# Insert two blank rows after the existing last employee row (25) - this
# pushes the trailing blank rows + footer rows down by 2 automatically.
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(27).Insert()

# The row that used to be the last (bold, bottom-bordered) employee row is
# still row 25. Copy that special formatting onto the new last row (27)
# before we overwrite row 25's formatting.
$ws.Range("B25:J25").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)

# Re-style the old last row (25) back to the normal (non-bold) row style,
# and give the same normal style to the other newly inserted row (26).
$ws.Range("B17:J17").Copy()
$ws.Range("B25:J26").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Write the final 12-row employee table (B16:G27)
# ---------------------------------------------------------------------------
$data = @(
  @("CC","1128056957","ERIK NAYID GARCIA LEDESMA","2104",10902,908526),
  @("CC","9145510","HERNAN MIGUEL SALCEDO VIVIERO","2104",6057,908526),
  @("CC","1047486409","WENDY JULIETH CARMONA MIRANDA","2104",23016,908526),
  @("CC","1143361451","ROBERTO CARLOS CARVAJAL CASTILLO","2104",10902,908526),
  @("CC","1050966362","ANGELLO ALEXANDER BUSTOS CABALLERO","2104",10902,908526),
  @("CC","73134020","JOSE HUMBERTO FIERRO LLAMAS","2104",2786,1044804),
  @("CC","1128058432","IVAN ANDRES DIAZ FERNANDEZ","2104",27861,908526),
  @("CC","1007170122","ALEXANDER MARRUGO JUNCO","2104",10902,908526),
  @("CC","1047510382","ANSELMO YANES MIRANDA","2104",4845,908526),
  @("CC","1127612794","LUIS ANSELMO YANES MIRANDA","2104",20593,908526),
  @("CC","1007981041","JOSE DAVID VERGARA OSPINA","2104",10902,908526),
  @("CC","1003061874","JOEL ANTONIO BULASCO CABRIA","2104",6057,908526)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Update the summary header fields
# ---------------------------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 145725
# Cant. Trabajadores
$ws.Range("C13").Value = 12

Write-Host "Employee table rewritten (12 rows) and header totals updated."
